# econ: anchor machines for wheeled harvester and forwarders
#
# Inserts two new parameter columns on the "parameterization" sheet:
#   - addOnWinchCableLength (350)   inserted immediately before ctlHaulHours (old col BF)
#   - anchorSMh (71.5)              inserted immediately before grappleYardingConstant
#                                    (old col CD, which after the first insert sits at CE)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("parameterization")
$ws.Activate()

# --- Insert the two new columns (real column inserts so every formula that
#     references a shifted cell gets its reference updated automatically) ---
$ws.Columns("BF:BF").Insert()
$ws.Columns("CE:CE").Insert()

# --- Populate "anchorSMh" first so it claims the lower shared-string index,
#     matching the order the strings were authored in the workbook ---
$ws.Range("CE1").Value = "anchorSMh"
$ws.Range("CE2").Value = 71.5
$ws.Range("CE2").NumberFormat = "0.00"

# --- Then populate "addOnWinchCableLength" ---
$ws.Range("BF1").Value = "addOnWinchCableLength"
$ws.Range("BF2").Value = 350
$ws.Range("BF2").NumberFormat = "0"

# --- Match the updated selection recorded in the sheet view ---
$ws.Range("BF1:BF2").Select()
